# Updates cryptos list prices / 1h volume percentages (GitHub Actions refresh).
# Cells whose new text would otherwise be auto-interpreted as a number
# (e.g. "605.44") are written with a leading apostrophe so they stay text,
# matching the original inlineStr cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.416.52'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.552.43'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''605.44'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '''144.59'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('D7').Value = '3.551.60'
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '''0.496'
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '''7.93'
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '4.154.30'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '''30.00'
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').Value = '3.561.50'
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('D17').Value = '66.459.47'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').Value = '''11.58'
$ws.Range('E19').Value = '  +6.06%  '
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').Value = '''14.83'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').Value = '''430.89'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').Value = '''0.609'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').Value = '3.691.89'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').Value = '''9.22'
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').Value = '''7.97'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '''1.45'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').Value = '3.547.04'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('D34').Value = '''25.32'
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').Value = '''0.153'
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('D36').Value = '''7.85'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('D39').Value = '''5.58'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').Value = '''174.72'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').Value = '''0.0848'
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').Value = '''5.18'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '''0.888'
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('D45').Value = '''46.11'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '''2.54'
$ws.Range('E47').Value = '  +5.89%  '
$ws.Range('D48').Value = '''1.18'
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('D49').Value = '''25.07'
$ws.Range('E49').Value = '  -3.78%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '''23.44'
$ws.Range('E51').Value = '  +4.08%  '
